$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(3.3, 3.4, 2.15, 2.4, 1.025, 0.9, 0.975, 0.85, 0.825, 1.1, 1.2, 1.3, 1.4, 1.5, 1.6, 1.7, 1.8, 1.9, 2.0, 2.1, 2.2, 2.3, 2.5, 2.6, 2.7, 2.8, 2.9, 3.0, 3.1, 3.2, 3.5)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(400 + $i, 2).Value2 = $vals[$i]
}
Write-Host "Done"
